# Add a header to the document containing the questionnaire number.
# (Questionnaire 40 -> "Questionnaire 40", centered, Arial 12pt, using the
# built-in "Header" paragraph style.)

$d = $word.ActiveDocument

# First (and only) section gets the new default header.
$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)   # wdHeaderFooterPrimary

$headerRange = $header.Range
$headerRange.InsertAfter("Questionnaire 40")

# Apply the "Header" paragraph style and center it.
$headerRange.Paragraphs.Item(1).Style = "Header"
$headerRange.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

# Format just the visible text run (exclude the trailing paragraph mark)
# so the run-level formatting lands on the <w:r>, not the paragraph mark.
$textRange = $header.Range.Duplicate
$textRange.End = $textRange.End - 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12
